# Applies the "Updated cryptos list" diff: refreshed Price/Volume(1h)
# figures, plus a WrappedEther/WrappedBTC row swap (rows 18-19).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.780.61"
$ws.Range("E2").Value = "  -2.35%  "
$ws.Range("D3").Value = "3.501.49"
$ws.Range("E3").Value = "  -3.31%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'586.01"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "'169.47"
$ws.Range("E6").Value = "  -3.41%  "
$ws.Range("D7").Value = "'0.608"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "3.494.85"
$ws.Range("E8").Value = "  -3.33%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  -4.65%  "
$ws.Range("D11").Value = "'6.76"
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("E12").Value = "  -4.92%  "
$ws.Range("D13").Value = "'46.79"
$ws.Range("E13").Value = "  -3.20%  "
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").Value = "4.070.45"
$ws.Range("E15").Value = "  -3.21%  "
$ws.Range("D16").Value = "'8.36"
$ws.Range("E16").Value = "  -6.19%  "
$ws.Range("D17").Value = "'609.96"
$ws.Range("E17").Value = "  -9.45%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "68.857.81"
$ws.Range("E18").Value = "  -2.27%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.503.10"
$ws.Range("E19").Value = "  -3.38%  "
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "'17.29"
$ws.Range("E21").Value = "  -2.59%  "
$ws.Range("D22").Value = "'11.06"
$ws.Range("E22").Value = "  -3.38%  "
$ws.Range("E23").Value = "  -6.31%  "
$ws.Range("D24").Value = "'15.61"
$ws.Range("E24").Value = "  -8.46%  "
$ws.Range("E25").Value = "  -3.86%  "
$ws.Range("D26").Value = "'3.81"
$ws.Range("E26").Value = "  -2.61%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").Value = "'2.59"
$ws.Range("E28").Value = "  -7.01%  "
$ws.Range("E29").Value = "  -6.82%  "
$ws.Range("D30").Value = "'32.33"
$ws.Range("E30").Value = "  -6.57%  "
$ws.Range("D31").Value = "'3.10"
$ws.Range("E31").Value = "  -5.63%  "
$ws.Range("D32").Value = "'8.40"
$ws.Range("E32").Value = "  -7.02%  "
$ws.Range("E33").Value = "  -5.72%  "
$ws.Range("E34").Value = "  -8.93%  "
$ws.Range("D35").Value = "'617.39"
$ws.Range("E35").Value = "  +6.38%  "
$ws.Range("E36").Value = "  -3.34%  "
$ws.Range("D37").Value = "'0.101"
$ws.Range("E37").Value = "  -4.96%  "
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("D39").Value = "'3.41"
$ws.Range("E39").Value = "  -14.21%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'0.0441"
$ws.Range("E41").Value = "  -2.28%  "
$ws.Range("E42").Value = "  -4.93%  "
$ws.Range("D43").Value = "3.367.92"
$ws.Range("E43").Value = "  -5.52%  "
$ws.Range("D44").Value = "'0.323"
$ws.Range("E44").Value = "  -5.93%  "
$ws.Range("D45").Value = "'32.50"
$ws.Range("E45").Value = "  -5.42%  "
$ws.Range("D46").Value = "0.0₃0687"
$ws.Range("E46").Value = "  -5.68%  "
$ws.Range("D47").Value = "'2.50"
$ws.Range("E47").Value = "  -6.61%  "
$ws.Range("E48").Value = "  -4.17%  "
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").Value = "'133.00"
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("E51").Value = "  +12.70%  "
